# Apply "Uploading new test values." edit to the Phase3 sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Phase3")

# Update the E column test values (and let Excel recalc F = LN(E))
$ws.Range("E7:E10").Value = 0.6
$ws.Range("E11:E14").Value = 0.7

# Update the selection shown in the sheet view
$ws.Activate()
$ws.Range("E21").Select()
